$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for columns D, L, M, N, O, P, Q, R, S, T for rows 2-11
# (A,B,C,E,F,G,H,I,J,K are unchanged across all rows)

$data = @{
    2  = @{ D = 44503; L = "Primera"; M = 50;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 2800; T = 10 }
    3  = @{ D = 44483; L = "Primera"; M = 35;  N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";   S = 2000; T = 5  }
    4  = @{ D = 44488; L = "Primera"; M = 100; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 5 kilos";  R = "La Ligua";                S = 2400; T = 5  }
    5  = @{ D = 44466; L = "Primera"; M = 80;  N = 11000; O = 11000; P = 11000; Q = "`$/bandeja 5 kilos";  R = "La Ligua";                S = 2200; T = 5  }
    6  = @{ D = 44166; L = "Segunda"; M = 20;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";    R = "La Ligua";                S = 667;  T = 18 }
    7  = @{ D = 44511; L = "Primera"; M = 45;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes";  S = 2800; T = 10 }
    8  = @{ D = 44511; L = "Primera"; M = 45;  N = 3200;  O = 3200;  P = 3200;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 320;  T = 10 }
    9  = @{ D = 44515; L = "Primera"; M = 80;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes";  S = 2800; T = 10 }
    10 = @{ D = 44519; L = "Primera"; M = 30;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 2800; T = 10 }
    11 = @{ D = 44496; L = "Primera"; M = 55;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 2800; T = 10 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T   # T: Kg / unidad
}
